$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "M3"
$ws.Range("E2").Value = "A1"
$ws.Range("H2").Value = "DO"
$ws.Range("I2").Value = "M1"
$ws.Range("J2").Value = "DO"
$ws.Range("K2").Value = "A1"
$ws.Range("N2").Value = "A1"
$ws.Range("O2").Value = "M3"
$ws.Range("R2").Value = "M1"
$ws.Range("V2").Value = "M1"
$ws.Range("W2").Value = "A1"
$ws.Range("X2").Value = "A1"
$ws.Range("Y2").Value = "DO"
$ws.Range("AB2").Value = "M3"
$ws.Range("AC2").Value = "M1"
$ws.Range("B3").Value = "DO"
$ws.Range("E3").Value = "M2"
$ws.Range("F3").Value = "A1"
$ws.Range("H3").Value = "A1"
$ws.Range("I3").Value = "DO"
$ws.Range("J3").Value = "A1"
$ws.Range("L3").Value = "M1"
$ws.Range("M3").Value = "M3"
$ws.Range("N3").Value = "M1"
$ws.Range("O3").Value = "A1"
$ws.Range("P3").Value = "M2"
$ws.Range("Q3").Value = "DO"
$ws.Range("R3").Value = "A1"
$ws.Range("S3").Value = "M1"
$ws.Range("T3").Value = "M3"
$ws.Range("U3").Value = "M1"
$ws.Range("V3").Value = "A2"
$ws.Range("X3").Value = "DO"
$ws.Range("Y3").Value = "A1"
$ws.Range("Z3").Value = "M1"
$ws.Range("AA3").Value = "M1"
$ws.Range("AB3").Value = "M2"
$ws.Range("AC3").Value = "A2"
$ws.Range("B4").Value = "A1"
$ws.Range("D4").Value = "M1"
$ws.Range("G4").Value = "M1"
$ws.Range("I4").Value = "A1"
$ws.Range("J4").Value = "DO"
$ws.Range("K4").Value = "M1"
$ws.Range("L4").Value = "A1"
$ws.Range("O4").Value = "M1"
$ws.Range("R4").Value = "M1"
$ws.Range("T4").Value = "M3"
$ws.Range("V4").Value = "M1"
$ws.Range("W4").Value = "A1"
$ws.Range("X4").Value = "M3"
$ws.Range("Y4").Value = "M1"
$ws.Range("AA4").Value = "DO"
$ws.Range("B5").Value = "DO"
$ws.Range("D5").Value = "M2"
$ws.Range("G5").Value = "M2"
$ws.Range("H5").Value = "M2"
$ws.Range("I5").Value = "DO"
$ws.Range("K5").Value = "M1"
$ws.Range("L5").Value = "M2"
$ws.Range("M5").Value = "M3"
$ws.Range("N5").Value = "M2"
$ws.Range("O5").Value = "M2"
$ws.Range("P5").Value = "DO"
$ws.Range("Q5").Value = "M1"
$ws.Range("R5").Value = "M2"
$ws.Range("S5").Value = "M1"
$ws.Range("T5").Value = "M1"
$ws.Range("U5").Value = "M1"
$ws.Range("V5").Value = "M2"
$ws.Range("W5").Value = "M1"
$ws.Range("Y5").Value = "M2"
$ws.Range("AB5").Value = "M3"
$ws.Range("AC5").Value = "M1"
$ws.Range("B6").Value = "A2"
$ws.Range("D6").Value = "M2"
$ws.Range("E6").Value = "A2"
$ws.Range("F6").Value = "M1"
$ws.Range("H6").Value = "DO"
$ws.Range("I6").Value = "A2"
$ws.Range("J6").Value = "M1"
$ws.Range("M6").Value = "M1"
$ws.Range("N6").Value = "M3"
$ws.Range("O6").Value = "DO"
$ws.Range("P6").Value = "M1"
$ws.Range("Q6").Value = "M1"
$ws.Range("S6").Value = "M1"
$ws.Range("U6").Value = "M1"
$ws.Range("V6").Value = "A2"
$ws.Range("W6").Value = "DO"
$ws.Range("X6").Value = "M1"
$ws.Range("Y6").Value = "M1"
$ws.Range("Z6").Value = "M2"
$ws.Range("AB6").Value = "M1"
$ws.Range("AC6").Value = "A1"
$ws.Range("B7").Value = "M3"
$ws.Range("C7").Value = "A1"
$ws.Range("D7").Value = "A1"
$ws.Range("E7").Value = "M1"
$ws.Range("F7").Value = "DO"
$ws.Range("H7").Value = "A1"
$ws.Range("I7").Value = "M3"
$ws.Range("J7").Value = "A1"
$ws.Range("K7").Value = "A1"
$ws.Range("L7").Value = "DO"
$ws.Range("M7").Value = "A1"
$ws.Range("N7").Value = "A1"
$ws.Range("O7").Value = "A1"
$ws.Range("Q7").Value = "A1"
$ws.Range("T7").Value = "DO"
$ws.Range("W7").Value = "M3"
$ws.Range("AA7").Value = "A1"
$ws.Range("AC7").Value = "DO"
$ws.Range("C8").Value = "A1"
$ws.Range("F8").Value = "DO"
$ws.Range("H8").Value = "A2"
$ws.Range("I8").Value = "A1"
$ws.Range("L8").Value = "DO"
$ws.Range("N8").Value = "M1"
$ws.Range("O8").Value = "A2"
$ws.Range("Q8").Value = "A2"
$ws.Range("R8").Value = "A2"
$ws.Range("T8").Value = "A2"
$ws.Range("U8").Value = "M2"
$ws.Range("V8").Value = "DO"
$ws.Range("W8").Value = "DO"
$ws.Range("AB8").Value = "M1"
$ws.Range("AC8").Value = "A2"
$ws.Range("B9").Value = "M1"
$ws.Range("C9").Value = "DO"
$ws.Range("F9").Value = "M3"
$ws.Range("H9").Value = "M1"
$ws.Range("I9").Value = "M2"
$ws.Range("J9").Value = "M2"
$ws.Range("K9").Value = "M1"
$ws.Range("L9").Value = "M2"
$ws.Range("M9").Value = "DO"
$ws.Range("O9").Value = "M1"
$ws.Range("P9").Value = "M2"
$ws.Range("Q9").Value = "M2"
$ws.Range("S9").Value = "DO"
$ws.Range("U9").Value = "A1"
$ws.Range("V9").Value = "M1"
$ws.Range("W9").Value = "M1"
$ws.Range("Y9").Value = "DO"
$ws.Range("AA9").Value = "M2"
$ws.Range("AB9").Value = "A2"
$ws.Range("AC9").Value = "M1"
$ws.Range("B10").Value = "M1"
$ws.Range("C10").Value = "M2"
$ws.Range("E10").Value = "M1"
$ws.Range("F10").Value = "A1"
$ws.Range("H10").Value = "DO"
$ws.Range("I10").Value = "M2"
$ws.Range("J10").Value = "M1"
$ws.Range("K10").Value = "PH"
$ws.Range("L10").Value = "A1"
$ws.Range("M10").Value = "M3"
$ws.Range("N10").Value = "M3"
$ws.Range("O10").Value = "DO"
$ws.Range("P10").Value = "A1"
$ws.Range("Q10").Value = "PH"
$ws.Range("R10").Value = "PH"
$ws.Range("S10").Value = "M1"
$ws.Range("U10").Value = "PH"
$ws.Range("V10").Value = "DO"
$ws.Range("W10").Value = "A2"
$ws.Range("X10").Value = "M2"
$ws.Range("Y10").Value = "M1"
$ws.Range("Z10").Value = "PH"
$ws.Range("AA10").Value = "M2"
$ws.Range("AB10").Value = "M2"
$ws.Range("AC10").Value = "DO"
